$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("area_mixre")
$ws1.Range("B2").Value = 207
$ws1.Range("B3").Value = 2.967776383623771
$ws1.Range("B4").Value = 2.783359747402721
$ws1.Range("B5").Value = 0.1270376448035983
$ws1.Range("B6").Value = 0.7454624506998567
$ws1.Range("B7").Value = 2.104469993336648
$ws1.Range("B8").Value = 4.534080824373094

$ws4 = $wb.Worksheets.Item("area_pop_sum")
$ws4.Range("B3").Value = 1128691
$ws4.Range("B4").Value = 1840.849159771859
